$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage (the source sheet
# stores every Coin/Link/Price/Volume cell as text, e.g. "60.928.73"
# or "1.00" must NOT be auto-converted to a number by COM type inference).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue 'D2' '60.928.73'
Set-TextValue 'E2' '  -0.73%  '

Set-TextValue 'D3' '3.392.14'
Set-TextValue 'E3' '  -1.16%  '

Set-TextValue 'D4' '1.00'

Set-TextValue 'D5' '570.23'
Set-TextValue 'E5' '  -0.88%  '

Set-TextValue 'D6' '142.11'
Set-TextValue 'E6' '  -2.17%  '

Set-TextValue 'D7' '3.392.84'
Set-TextValue 'E7' '  -1.08%  '

Set-TextValue 'E8' '  -0.03%  '

Set-TextValue 'D9' '0.476'
Set-TextValue 'E9' '  -0.57%  '

Set-TextValue 'D10' '7.54'
Set-TextValue 'E10' '  -1.73%  '

Set-TextValue 'E11' '  -1.51%  '

Set-TextValue 'D13' '3.971.27'
Set-TextValue 'E13' '  -1.28%  '

Set-TextValue 'D14' '28.34'
Set-TextValue 'E14' '  +1.42%  '

Set-TextValue 'E15' '  +1.94%  '

Set-TextValue 'E16' '  -1.21%  '

Set-TextValue 'D17' '3.395.07'
Set-TextValue 'E17' '  -1.01%  '

Set-TextValue 'D18' '60.999.26'
Set-TextValue 'E18' '  -0.81%  '

Set-TextValue 'D19' '6.18'
Set-TextValue 'E19' '  -1.56%  '

Set-TextValue 'E20' '  -2.11%  '

Set-TextValue 'E21' '  -4.45%  '

Set-TextValue 'D22' '385.26'
Set-TextValue 'E22' '  -2.49%  '

Set-TextValue 'E23' '  -1.28%  '

Set-TextValue 'D24' '74.15'
Set-TextValue 'E24' '  +1.00%  '

Set-TextValue 'E25' '  +0.39%  '

Set-TextValue 'E26' '  -4.63%  '

Set-TextValue 'D27' '3.531.87'
Set-TextValue 'E27' '  -1.12%  '

Set-TextValue 'E28' '  -0.58%  '

Set-TextValue 'D29' '0.999'
Set-TextValue 'E29' '  -0.10%  '

Set-TextValue 'E30' '  -2.81%  '

Set-TextValue 'E31' '  -2.90%  '

Set-TextValue 'B32' 'Fetch.AI'
Set-TextValue 'C32' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D32' '1.43'
Set-TextValue 'E32' '  -2.30%  '

Set-TextValue 'B33' 'PancakeSwap'
Set-TextValue 'C33' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D33' '2.13'
Set-TextValue 'E33' '  -1.83%  '

Set-TextValue 'E34' '  -0.03%  '

Set-TextValue 'D35' '23.54'
Set-TextValue 'E35' '  -1.88%  '

Set-TextValue 'D36' '6.98'
Set-TextValue 'E36' '  -0.43%  '

Set-TextValue 'D37' '167.65'
Set-TextValue 'E37' '  +0.03%  '

Set-TextValue 'D38' '3.422.41'
Set-TextValue 'E38' '  -1.14%  '

Set-TextValue 'D39' '4.98'
Set-TextValue 'E39' '  -2.52%  '

Set-TextValue 'E40' '  -4.80%  '

Set-TextValue 'D41' '27.79'
Set-TextValue 'E41' '  +2.88%  '

Set-TextValue 'E42' '  -1.05%  '

Set-TextValue 'E43' '  +0.02%  '

Set-TextValue 'E44' '  -2.38%  '

Set-TextValue 'E45' '  +0.04%  '

Set-TextValue 'E46' '  -1.27%  '

Set-TextValue 'E47' '  -3.43%  '

Set-TextValue 'D48' '1.13'
Set-TextValue 'E48' '  -1.97%  '

Set-TextValue 'D49' '2.488.41'
Set-TextValue 'E49' '  -3.95%  '

Set-TextValue 'E50' '  -1.47%  '

Set-TextValue 'D51' '23.04'
Set-TextValue 'E51' '  -0.25%  '
